$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. The stray "_GoBack" bookmark currently sits right after
#    "Tried addressing this." -- drop it; Word will re-plant a fresh
#    one at whatever the most-recently-edited spot turns out to be.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Swap the two runs "Maybe reframe associations as first order
#    associations" + "?" for a single new sentence.
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(
    "Maybe reframe associations as first order associations?",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Not sure the best way to approach this", 2) | Out-Null

# $rng now covers exactly the freshly-inserted replacement text.
$editEnd = $rng.End

# ------------------------------------------------------------------
# 3. Re-plant "_GoBack" immediately after the new sentence (last
#    thing in that paragraph, same as Word does after any edit).
#    A zero-length range sitting exactly on the paragraph-mark slot
#    doesn't round-trip cleanly, so nudge past it with a throwaway
#    marker, bookmark the point before the marker, then remove it.
# ------------------------------------------------------------------
$marker = $d.Range($editEnd, $editEnd)
$marker.InsertAfter("@@bm@@")

$bmPoint = $d.Range($editEnd, $editEnd)
$d.Bookmarks.Add("_GoBack", $bmPoint)

$d.Range($editEnd, $editEnd + 6).Delete()
